# update stategy setting for czce night trading
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("策略更新")

# F4: fix the CZCE night-trading EV-file remark (drop the stray "f50" token)
$ws.Range("F4").Value = "更新夜盘ev文件"

# D4: move the CZCE night-trading date forward by one day (2017-01-09 -> 2017-01-10)
$ws.Range("D4").Value = 42745

# Row 5 used to hold the "上期 / 夜盘 / 调整一下夜盘的手数..." note, which is no
# longer relevant - clear it out. B5:D5 keep their existing formatting, while
# E5:F5 are fully cleared (content + formatting) so nothing is left behind.
$ws.Range("B5:D5").ClearContents()
$ws.Range("E5:F5").Clear()
$ws.Rows.Item(5).AutoFit()

# Restore selection to E4 (matches the refreshed view after the edits above).
$ws.Range("E4").Select()
